# Applies the "updating surface area and scripts" edit:
#  - Adds a new "Level" shared string
#  - Duplicates the Run1/Run2 PI-curve tables (rows 2-14 and 15-24) into a
#    second copy starting at row 27, this time paired up with a new
#    "Level" (L) column that averages the two runs together row-for-row,
#    and a simple running index in column M.
#  - Leaves rows 1-24 completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 27: header row (same headers as row 1, plus a new "Level" column)
# ---------------------------------------------------------------------
$headers = @("percent","chamber1","chamber2","chamber3","chamber4","chamber5","chamber6","chamber7","chamber8","average","Run","Level")
for ($c = 0; $c -lt $headers.Length; $c++) {
  $ws.Cells.Item(27, $c + 1).Value = $headers[$c]
}

# ---------------------------------------------------------------------
# Data blocks: (row, percentLabel, chamber1..8 values)
# ---------------------------------------------------------------------
$block1 = @(
  @(28, "0",   @(0,0,0,0,0,0,0,0)),
  @(29, "10",  @(33,34,45,46,40,36,33,31)),
  @(30, "20",  @(85,91,119,117,101,96,87,78)),
  @(31, "30",  @(127,134,176,175,145,140,130,119)),
  @(32, "40",  @(165,172,219,234,202,188,172,151)),
  @(33, "60",  @(232,258,325,333,286,267,238,224)),
  @(34, "80",  @(292,314,414,407,355,329,295,281)),
  @(35, "90",  @(320,339,442,444,368,354,315,295)),
  @(36, "100", @(333,367,465,468,418,378,328,317)),
  @(37, "100+zip", @(505,485,600,529,611,455,478,380))
)

$block2 = @(
  @(38, "0",   @(0,0,0,0,0,0,0,0)),
  @(39, "10",  @(29,30,33,32,36,30,35,28)),
  @(40, "20",  @(73,77,85,83,93,80,88,73)),
  @(41, "30",  @(112,119,130,127,140,120,130,112)),
  @(42, "40",  @(150,157,174,171,185,160,174,148)),
  @(43, "60",  @(211,223,241,243,263,230,245,205)),
  @(44, "80",  @(262,279,300,298,326,285,304,260)),
  @(45, "125", @(311,325,350,349,385,335,365,309)),
  @(46, "125+3clicks", @(377,416,430,430,462,436,404,352)),
  @(47, "125+6clicks", @(484,492,553,545,593,565,485,418))
)

$numericLabels = @("0","10","20","30","40","50","60","70","80","90","100","125")

function Set-PercentCell($ws, $row, $label) {
  if ($numericLabels -contains $label) {
    $ws.Cells.Item($row, 1).Value = [double]$label
  } else {
    $ws.Cells.Item($row, 1).Value = $label
  }
}

# ---------------------------------------------------------------------
# Block 1 (rows 28-37): "Run1" data, column K = Run1
# ---------------------------------------------------------------------
foreach ($entry in $block1) {
  $row = $entry[0]
  $label = $entry[1]
  $vals = $entry[2]

  Set-PercentCell $ws $row $label

  for ($c = 0; $c -lt $vals.Length; $c++) {
    $ws.Cells.Item($row, 2 + $c).Value = $vals[$c]
  }

  $ws.Cells.Item($row, 10).Formula = "=AVERAGE(B$($row):I$($row))"
  $ws.Cells.Item($row, 10).NumberFormat = "0"

  $ws.Cells.Item($row, 11).Value = "Run1"
}

# ---------------------------------------------------------------------
# Block 2 (rows 38-47): "Run2" data, column K = Run2, plus Level (L)
# and running index (M) that pairs each row with the matching row in
# block 1 (10 rows above).
# ---------------------------------------------------------------------
for ($i = 0; $i -lt $block2.Length; $i++) {
  $entry = $block2[$i]
  $row = $entry[0]
  $label = $entry[1]
  $vals = $entry[2]
  $pairRow = $block1[$i][0]

  Set-PercentCell $ws $row $label

  for ($c = 0; $c -lt $vals.Length; $c++) {
    $ws.Cells.Item($row, 2 + $c).Value = $vals[$c]
  }

  $ws.Cells.Item($row, 10).Formula = "=AVERAGE(B$($row):I$($row))"
  $ws.Cells.Item($row, 10).NumberFormat = "0"

  $ws.Cells.Item($row, 11).Value = "Run2"

  $ws.Cells.Item($row, 12).Formula = "=AVERAGE(B$($pairRow):I$($pairRow),B$($row):I$($row))"
  $ws.Cells.Item($row, 12).NumberFormat = "0"

  $ws.Cells.Item($row, 13).Value = $i + 1
}

# ---------------------------------------------------------------------
# Final view state: scroll to show the new table and select L38:L47
# ---------------------------------------------------------------------
$ws.Range("L38:L47").Select()
try {
  $excel.ActiveWindow.ScrollRow = 23
  $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
